$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the now-unused trailing rows (45-51) so the sheet ends at row 44
$ws.Rows("45:51").Delete()

# Update the Firm / Time / Lawyers Registered data for rows 2-44
$ws.Cells.Item(2, 1).Value = "AL Goodbody"
$ws.Cells.Item(2, 2).Value = "40s"
$ws.Cells.Item(2, 3).Value = "2"
$ws.Cells.Item(3, 1).Value = "Dottir"
$ws.Cells.Item(3, 2).Value = "04s"
$ws.Cells.Item(3, 3).Value = "0"
$ws.Cells.Item(4, 1).Value = "RDJ"
$ws.Cells.Item(4, 2).Value = "04s"
$ws.Cells.Item(4, 3).Value = "1"
$ws.Cells.Item(5, 1).Value = "Gomez Acebo And Pombo"
$ws.Cells.Item(5, 2).Value = "08s"
$ws.Cells.Item(5, 3).Value = "2"
$ws.Cells.Item(6, 1).Value = "Frontier"
$ws.Cells.Item(6, 2).Value = "03min 26s"
$ws.Cells.Item(6, 3).Value = "2"
$ws.Cells.Item(7, 1).Value = "BonelliErede"
$ws.Cells.Item(7, 2).Value = "15s"
$ws.Cells.Item(7, 3).Value = "1"
$ws.Cells.Item(8, 1).Value = "HNA"
$ws.Cells.Item(8, 2).Value = "10s"
$ws.Cells.Item(8, 3).Value = "0"
$ws.Cells.Item(9, 1).Value = "DR And AJU"
$ws.Cells.Item(9, 2).Value = "06s"
$ws.Cells.Item(9, 3).Value = "1"
$ws.Cells.Item(10, 1).Value = "BCF Law"
$ws.Cells.Item(10, 2).Value = "16s"
$ws.Cells.Item(10, 3).Value = "1"
$ws.Cells.Item(11, 1).Value = "MSP"
$ws.Cells.Item(11, 2).Value = "23s"
$ws.Cells.Item(11, 3).Value = "0"
$ws.Cells.Item(12, 1).Value = "Nelligan Law"
$ws.Cells.Item(12, 2).Value = "26s"
$ws.Cells.Item(12, 3).Value = "1"
$ws.Cells.Item(13, 1).Value = "Beauchamps"
$ws.Cells.Item(13, 2).Value = "49s"
$ws.Cells.Item(13, 3).Value = "1"
$ws.Cells.Item(14, 1).Value = "Tompkins Wake"
$ws.Cells.Item(14, 2).Value = "35s"
$ws.Cells.Item(14, 3).Value = "1"
$ws.Cells.Item(15, 1).Value = "Reliance Corporate Advisors"
$ws.Cells.Item(15, 2).Value = "21s"
$ws.Cells.Item(15, 3).Value = "1"
$ws.Cells.Item(16, 1).Value = "Tiruchelvam Associates"
$ws.Cells.Item(16, 2).Value = "07s"
$ws.Cells.Item(16, 3).Value = "1"
$ws.Cells.Item(17, 1).Value = "Kinstellar"
$ws.Cells.Item(17, 2).Value = "37s"
$ws.Cells.Item(17, 3).Value = "3"
$ws.Cells.Item(18, 1).Value = "Matheson"
$ws.Cells.Item(18, 2).Value = "10s"
$ws.Cells.Item(18, 3).Value = "1"
$ws.Cells.Item(19, 1).Value = "Pedersoli"
$ws.Cells.Item(19, 2).Value = "01min 13s"
$ws.Cells.Item(19, 3).Value = "0"
$ws.Cells.Item(20, 1).Value = "Andersen"
$ws.Cells.Item(20, 2).Value = "06s"
$ws.Cells.Item(20, 3).Value = "1"
$ws.Cells.Item(21, 1).Value = "Arnesen IP"
$ws.Cells.Item(21, 2).Value = "15s"
$ws.Cells.Item(21, 3).Value = "0"
$ws.Cells.Item(22, 1).Value = "Oyen Wiggs"
$ws.Cells.Item(22, 2).Value = "05s"
$ws.Cells.Item(22, 3).Value = "1"
$ws.Cells.Item(23, 1).Value = "DCC Law"
$ws.Cells.Item(23, 2).Value = "39s"
$ws.Cells.Item(23, 3).Value = "2"
$ws.Cells.Item(24, 1).Value = "Zamfirescu Racoti Predoiu"
$ws.Cells.Item(24, 2).Value = "07s"
$ws.Cells.Item(24, 3).Value = "1"
$ws.Cells.Item(25, 1).Value = "Stephenson Harwood"
$ws.Cells.Item(25, 2).Value = "28s"
$ws.Cells.Item(25, 3).Value = "3"
$ws.Cells.Item(26, 1).Value = "Magnusson Law"
$ws.Cells.Item(26, 2).Value = "31s"
$ws.Cells.Item(26, 3).Value = "2"
$ws.Cells.Item(27, 1).Value = "Lee And Ko"
$ws.Cells.Item(27, 2).Value = "18s"
$ws.Cells.Item(27, 3).Value = "1"
$ws.Cells.Item(28, 1).Value = "StewartMcKelvey"
$ws.Cells.Item(28, 2).Value = "16s"
$ws.Cells.Item(28, 3).Value = "1"
$ws.Cells.Item(29, 1).Value = "Wildeboer Dellelce"
$ws.Cells.Item(29, 2).Value = "05s"
$ws.Cells.Item(29, 3).Value = "1"
$ws.Cells.Item(30, 1).Value = "DBHLaw"
$ws.Cells.Item(30, 2).Value = "12s"
$ws.Cells.Item(30, 3).Value = "1"
$ws.Cells.Item(31, 1).Value = "JGSA"
$ws.Cells.Item(31, 2).Value = "10s"
$ws.Cells.Item(31, 3).Value = "1"
$ws.Cells.Item(32, 1).Value = "Liedekerke"
$ws.Cells.Item(32, 2).Value = "38s"
$ws.Cells.Item(32, 3).Value = "1"
$ws.Cells.Item(33, 1).Value = "BNT"
$ws.Cells.Item(33, 2).Value = "24s"
$ws.Cells.Item(33, 3).Value = "3"
$ws.Cells.Item(34, 1).Value = "Barriston Law"
$ws.Cells.Item(34, 2).Value = "05s"
$ws.Cells.Item(34, 3).Value = "1"
$ws.Cells.Item(35, 1).Value = "Borenius"
$ws.Cells.Item(35, 2).Value = "17s"
$ws.Cells.Item(35, 3).Value = "1"
$ws.Cells.Item(36, 1).Value = "Ogletree Deakins"
$ws.Cells.Item(36, 2).Value = "13s"
$ws.Cells.Item(36, 3).Value = "1"
$ws.Cells.Item(37, 1).Value = "James And Wells"
$ws.Cells.Item(37, 2).Value = "15s"
$ws.Cells.Item(37, 3).Value = "2"
$ws.Cells.Item(38, 1).Value = "Brigrard Urrutia"
$ws.Cells.Item(38, 2).Value = "18s"
$ws.Cells.Item(38, 3).Value = "1"
$ws.Cells.Item(39, 1).Value = "Covenant Chambers"
$ws.Cells.Item(39, 2).Value = "19s"
$ws.Cells.Item(39, 3).Value = "1"
$ws.Cells.Item(40, 1).Value = "Dahl Law"
$ws.Cells.Item(40, 2).Value = "12s"
$ws.Cells.Item(40, 3).Value = "1"
$ws.Cells.Item(41, 1).Value = "Barnea And Co"
$ws.Cells.Item(41, 2).Value = "27s"
$ws.Cells.Item(41, 3).Value = "1"
$ws.Cells.Item(42, 1).Value = "Watson Farley And Williams"
$ws.Cells.Item(42, 2).Value = "11s"
$ws.Cells.Item(42, 3).Value = "3"
$ws.Cells.Item(43, 1).Value = "Cassidy Levy Kent"
$ws.Cells.Item(43, 2).Value = "22s"
$ws.Cells.Item(43, 3).Value = "2"
$ws.Cells.Item(44, 1).Value = "Conyers"
$ws.Cells.Item(44, 2).Value = "10s"
$ws.Cells.Item(44, 3).Value = "3"
